# Modelo de Análise dos Eventos.xlsx - edit matching commit:
# "Alterei o meu evento, deixei entre colchetes para ver se todos concordam"
#
# The sheet "Plan3" (3rd sheet, the active tab) has its D3 and D4 texts
# updated, and the current selection moves from D11 to D4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan3")

$ws.Range("D3").Value = "Motoboy passa as informações sobre a manutenção realizada no veículo"
$ws.Range("D4").Value = "O Gerente de Manutenção [disponibiliza as manutenções realizadas para o Administrativo] encaminha as informações para o Administrativo  no final da semana"

$ws.Activate()
$ws.Range("D4").Select()
